$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reassign feature rows 4 and 13 (Advanced Settings Menu, Iteration Tab) to Sebastian
$ws.Range("C4").Value = "Sebastian"
$ws.Range("C13").Value = "Sebastian"

# Update the active selection to C13, matching the author's last edit position
$ws.Range("C13").Select()
